$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

$values = @(
    "26, 64, 56, 34, 43",
    "36, 10, 39, 44, 3",
    "86, 55, 79, 13, 31",
    "74, 36, 83, 71, 23",
    "41, 74, 62, ",
    "78, 59, 34",
    "34, 11",
    "78, 59, 62, 2",
    "29, 34, 91",
    "34, 73",
    "34, 57",
    "41, 74, 84",
    "76, 70, 34",
    "34, 67",
    "34, 92, 50, 35",
    "21, 53, 34",
    "82, 50, 35",
    "34, 53, 35, 63",
    "60, 39",
    "45, 39, 54",
    "41, 48, 7, 15, 81",
    "46, 79, 13, 85, 14",
    "38, 51, 46, 79, 13",
    "38, 58",
    "42, 25",
    "66, 1, 71, 37",
    "18, 47, 85, 33",
    "28, 48, 6, 50, 24",
    "47, 85, 28, 49",
    "47, 18, 85",
    "90, 89",
    "34, 80, 39",
    "34, 72, 65",
    "48, 34, 69",
    "87, 39, 50",
    "75, 40",
    "34, 9",
    "88, 61",
    "86, 68, 77",
    "17, 16",
    "93, 32",
    "27, 22",
    "52, 19",
    "80, 66, 5, 34",
    "4, 39",
    "34, 8",
    "12, 79, 13"
)

$row = 5
foreach ($val in $values) {
    $ws.Cells.Item($row, 3).Value = $val
    $row = $row + 1
}

$ws.Activate()
$ws.Range("C53").Select()
